$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "bay" status cells in column E for rows 3-5, mirroring column N.
# Copy the (border-only, non-date) format from the matching N-column cells so the
# now-unused date-format style on E3/E4 gets dropped, then set the values.
$ws.Range("N3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "free"

$ws.Range("N4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "ok"

$ws.Range("E5").Value = "wait"

# Update the active selection on the sheet
$ws.Range("K18").Select()
